$wb = $excel.ActiveWorkbook

# --- ALC (sheet1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 778.5
$ws.Range("I33").Value = 803.7692
$ws.Range("K33").Value = 803.7692
$ws.Range("M33").Value = -574.7692
$ws.Range("H47").Value = 24999
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H70").Value = 6828.722
$ws.Range("I70").Value = 6537.5
$ws.Range("J70").Value = 6974.3335
$ws.Range("K70").Value = 19612.5
$ws.Range("L70").Value = 20923.0005
$ws.Range("M70").Value = -19342.5
$ws.Range("N70").Value = -21463.0005
$ws.Range("H73").Value = 6828.722
$ws.Range("I73").Value = 6537.5
$ws.Range("J73").Value = 6974.3335
$ws.Range("K73").Value = 19612.5
$ws.Range("L73").Value = 20923.0005
$ws.Range("M73").Value = -18676.5
$ws.Range("N73").Value = -22795.0005
$ws.Range("H76").Value = 5499.8
$ws.Range("I76").Value = 5499.7856
$ws.Range("K76").Value = 5499.7856
$ws.Range("M76").Value = -5184.7856
$ws.Range("H79").Value = 5499.8
$ws.Range("I79").Value = 5499.7856
$ws.Range("K79").Value = 5499.7856
$ws.Range("M79").Value = -4407.7856
$ws.Range("H88").Value = 2630.3333
$ws.Range("I88").Value = 875
$ws.Range("J88").Value = 2900.3845
$ws.Range("K88").Value = 875
$ws.Range("L88").Value = 2900.3845
$ws.Range("M88").Value = -469
$ws.Range("N88").Value = -3712.3845
$ws.Range("H91").Value = 2630.3333
$ws.Range("I91").Value = 875
$ws.Range("J91").Value = 2900.3845
$ws.Range("K91").Value = 875
$ws.Range("L91").Value = 2900.3845
$ws.Range("M91").Value = 529
$ws.Range("N91").Value = -5708.3845
$ws.Range("H112").Value = 1122.5358
$ws.Range("J112").Value = 1110.8077
$ws.Range("L112").Value = 3332.4231
$ws.Range("N112").Value = -5548.4231
$ws.Range("H132").Value = 1369.6666
$ws.Range("I132").Value = 1320.2858
$ws.Range("K132").Value = 3960.8574
$ws.Range("M132").Value = -1430.8574
$ws.Range("H135").Value = 2348.5715
$ws.Range("I135").Value = 1297.8
$ws.Range("K135").Value = 11680.2
$ws.Range("M135").Value = -9145.199999999999

# --- ARM (sheet2) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3665.9756
$ws.Range("I32").Value = 3187.4358
$ws.Range("J32").Value = 12997.5
$ws.Range("K32").Value = 3187.4358
$ws.Range("L32").Value = 12997.5
$ws.Range("M32").Value = -2900.4358
$ws.Range("N32").Value = -13571.5
$ws.Range("H102").Value = 3061
$ws.Range("I102").Value = 3061
$ws.Range("K102").Value = 3061
$ws.Range("M102").Value = -1439
$ws.Range("H132").Value = 1286.5
$ws.Range("J132").Value = 865
$ws.Range("L132").Value = 2595
$ws.Range("N132").Value = -7655

# --- BSM (sheet3) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1798.3914
$ws.Range("J20").Value = 1693.3334
$ws.Range("L20").Value = 1693.3334
$ws.Range("N20").Value = -2187.3334
$ws.Range("H86").Value = 1320.5714
$ws.Range("I86").Value = 1207.5
$ws.Range("K86").Value = 1207.5
$ws.Range("M86").Value = -84.5
$ws.Range("H89").Value = 1320.5714
$ws.Range("I89").Value = 1207.5
$ws.Range("K89").Value = 6037.5
$ws.Range("M89").Value = -421.5
$ws.Range("H107").Value = 3429.2856
$ws.Range("I107").Value = 2398.25
$ws.Range("K107").Value = 2398.25
$ws.Range("M107").Value = -478.25
$ws.Range("H134").Value = 1040.8
$ws.Range("I134").Value = 1001
$ws.Range("K134").Value = 3003
$ws.Range("M134").Value = -468

# --- CRP (sheet4) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 996.75
$ws.Range("I58").Value = 995.6667
$ws.Range("K58").Value = 995.6667
$ws.Range("M58").Value = -792.6667
$ws.Range("H132").Value = 1762.7858
$ws.Range("I132").Value = 1562.2727
$ws.Range("J132").Value = 2498
$ws.Range("K132").Value = 4686.8181
$ws.Range("L132").Value = 7494
$ws.Range("M132").Value = -2156.8181
$ws.Range("N132").Value = -12554
$ws.Range("H136").Value = 996.75
$ws.Range("I136").Value = 995.6667
$ws.Range("K136").Value = 2987.0001
$ws.Range("M136").Value = -437.0001000000002
$ws.Range("H141").Value = 280106.1
$ws.Range("J141").Value = 280106.1
$ws.Range("L141").Value = 280106.1
$ws.Range("N141").Value = -290466.1

# --- CUL (sheet5) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 22001308
$ws.Range("I4").Value = 22001308
$ws.Range("K4").Value = 66003924
$ws.Range("M4").Value = -66003812
$ws.Range("H18").Value = 923.8
$ws.Range("I18").Value = 923.8
$ws.Range("K18").Value = 2771.4
$ws.Range("M18").Value = -2602.4
$ws.Range("H70").Value = 12366.5
$ws.Range("I70").Value = 2466.5
$ws.Range("J70").Value = 15666.5
$ws.Range("K70").Value = 7399.5
$ws.Range("L70").Value = 46999.5
$ws.Range("M70").Value = -7084.5
$ws.Range("N70").Value = -47629.5
$ws.Range("H73").Value = 12366.5
$ws.Range("I73").Value = 2466.5
$ws.Range("J73").Value = 15666.5
$ws.Range("K73").Value = 7399.5
$ws.Range("L73").Value = 46999.5
$ws.Range("M73").Value = -6307.5
$ws.Range("N73").Value = -49183.5
$ws.Range("H75").Value = 3278.1428
$ws.Range("I75").Value = 3474.5
$ws.Range("J75").Value = 3199.6
$ws.Range("K75").Value = 10423.5
$ws.Range("L75").Value = 9598.799999999999
$ws.Range("M75").Value = -9425.5
$ws.Range("N75").Value = -11594.8
$ws.Range("H78").Value = 3278.1428
$ws.Range("I78").Value = 3474.5
$ws.Range("J78").Value = 3199.6
$ws.Range("K78").Value = 31270.5
$ws.Range("L78").Value = 28796.4
$ws.Range("M78").Value = -26278.5
$ws.Range("N78").Value = -38780.39999999999
$ws.Range("H134").Value = 15570
$ws.Range("I134").Value = 6710
$ws.Range("K134").Value = 20130
$ws.Range("M134").Value = -15060
$ws.Range("H141").Value = 13874.5
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- GSM (sheet6) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5252.625
$ws.Range("I70").Value = 4872.25
$ws.Range("J70").Value = 5633
$ws.Range("K70").Value = 4872.25
$ws.Range("L70").Value = 5633
$ws.Range("M70").Value = -4602.25
$ws.Range("N70").Value = -6173
$ws.Range("H73").Value = 5252.625
$ws.Range("I73").Value = 4872.25
$ws.Range("J73").Value = 5633
$ws.Range("K73").Value = 4872.25
$ws.Range("L73").Value = 5633
$ws.Range("M73").Value = -3936.25
$ws.Range("N73").Value = -7505
$ws.Range("H107").Value = 416.18182
$ws.Range("J107").Value = 481.66666
$ws.Range("L107").Value = 481.66666
$ws.Range("N107").Value = -4321.66666
$ws.Range("H126").Value = 2999
$ws.Range("I126").Value = 2999
$ws.Range("K126").Value = 8997
$ws.Range("M126").Value = -6527
$ws.Range("H132").Value = 1249
$ws.Range("I132").Value = 1249
$ws.Range("K132").Value = 3747
$ws.Range("M132").Value = -1217

# --- LTW (sheet7) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4266.4287
$ws.Range("I22").Value = 4769
$ws.Range("K22").Value = 4769
$ws.Range("M22").Value = -4474
$ws.Range("H27").Value = 4266.4287
$ws.Range("I27").Value = 4769
$ws.Range("K27").Value = 4769
$ws.Range("M27").Value = -4662
$ws.Range("H132").Value = 4182
$ws.Range("I132").Value = 3977.75
$ws.Range("K132").Value = 11933.25
$ws.Range("M132").Value = -9403.25
$ws.Range("H136").Value = 8652.1
$ws.Range("I136").Value = 10610.143
$ws.Range("J136").Value = 4083.3333
$ws.Range("K136").Value = 31830.429
$ws.Range("L136").Value = 12249.9999
$ws.Range("M136").Value = -29280.429
$ws.Range("N136").Value = -17349.9999

# --- WVR (sheet8) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3864.111
$ws.Range("I132").Value = 2472.25
$ws.Range("J132").Value = 14999
$ws.Range("K132").Value = 7416.75
$ws.Range("L132").Value = 44997
$ws.Range("M132").Value = -4886.75
$ws.Range("N132").Value = -50057
$ws.Range("H136").Value = 2800.923
$ws.Range("I136").Value = 2884.3333
$ws.Range("J136").Value = 1800
$ws.Range("K136").Value = 8652.999899999999
$ws.Range("L136").Value = 5400
$ws.Range("M136").Value = -6102.999899999999
$ws.Range("N136").Value = -10500
